$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing birthday date for row 2 (John) ---
$ws.Range("C2").Value = 32970

# --- Add new row 3 data (Alina) ---
$ws.Range("A3").Value = 1904
$ws.Range("B3").Value = "Alina"

# Birthday date for Alina - copy formatting from C2 (date style) then set value
$ws.Range("C2").Copy($ws.Range("C3")) | Out-Null
$ws.Range("C3").Value = 36623

# Email (with hyperlink) for Alina - copy formatting from D2 (hyperlink style),
# add the hyperlink, then restore the exact formatting/value afterwards so the
# resulting cell style matches the existing hyperlink style used by D2.
$ws.Range("D2").Copy($ws.Range("D3")) | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:akis441981@gmail.com") | Out-Null
$ws.Range("D2").Copy($ws.Range("D3")) | Out-Null
$ws.Range("D3").Value = "akis441981@gmail.com"

$ws.Range("E3").Value = 87465487887
$ws.Range("F3").Value = "fknvbckjbnkcjFHFHGFkhkjhkjhkj76576jdfkjdh"

# Match the saved selection state from the edited workbook
$ws.Range("C3").Select() | Out-Null
